$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Granular"
$ws.Range("E2").Value = $false
$ws.Range("E3").Value = $true
$ws.Range("E4").Value = $false
$ws.Range("E5").Value = $true
$ws.Range("E6").Value = $true

$ws.Range("E3").Select()
